$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell K1, matching the style of the existing header row (A1:J1)
$ws.Range("K1").Value = "intervention_type"
$ws.Range("K1").Style = $ws.Range("J1").Style

# Add new data cell K2
$ws.Range("K2").Value = "OTHER"
